# Correct "TDY1770" to "TDY1779" in the strain column (column E)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
foreach ($cell in $used.Cells) {
    if ($cell.Value2 -eq "TDY1770") {
        $cell.Value2 = "TDY1779"
    }
}
